# Assignment 3: Coding Basics updates
#
# - Data sheet, A4: was the text "sixty" -> becomes the number 60
# - Data sheet, B11: was 7000 -> becomes 70
# - The now-unused shared string "sixty" is dropped automatically once
#   nothing references it any more.
# - "Data" becomes the active/selected sheet (cell G10 selected) instead
#   of "Codebook" (which keeps its own last selection of C10).

$wb = $excel.ActiveWorkbook

$wsData     = $wb.Worksheets.Item(1)
$wsCodebook = $wb.Worksheets.Item(2)

# Fix up the two data values.
$wsData.Range("A4").Value = 60
$wsData.Range("B11").Value = 70

# Record Codebook's (unchanged) selection before switching away from it.
[void]$wsCodebook.Range("C10").Select()

# Switch to the Data sheet and select G10, leaving Data as the active tab.
[void]$wsData.Activate()
[void]$wsData.Range("G10").Select()
